# This workbook is a weekly price log. A new week's worth of Cereza
# ("Rainier" variety) price rows is inserted right after the existing
# row 73, which pushes all the following data rows down by 4 positions
# (old rows 74:146 become 78:150). The sheet's used range therefore
# grows from A1:T146 to A1:T150.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows at position 74; everything from row 74 down
# (through the former last row, 146) shifts down to make room.
$ws.Rows("74:77").Insert()

# Values that differ between the 4 new rows: Calidad, Volumen,
# Precio minimo/maximo/promedio ponderado, and Precio $/Kg.
$newRowNums = @(74, 75, 76, 77)
$quality    = @("Primera", "Segunda", "Primera", "Segunda")
$volume     = @(120,        60,        120,       60)
$minPrice   = @(5000,       4000,      5000,      4000)
$maxPrice   = @(6000,       4000,      6000,      4000)
$avgPrice   = @(5500,       4000,      5500,      4000)
$perKg      = @(550,        400,       550,       400)

for ($i = 0; $i -lt $newRowNums.Length; $i++) {
    $r = $newRowNums[$i]

    $ws.Cells.Item($r, 1).Value  = 7
    $ws.Cells.Item($r, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
    $ws.Cells.Item($r, 3).Value  = "Ñuble"

    $ws.Cells.Item($r, 4).Value        = [DateTime]"2022-12-28"
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($r, 5).Value  = 16
    $ws.Cells.Item($r, 6).Value  = "Fruta"
    $ws.Cells.Item($r, 7).Value  = 100103
    $ws.Cells.Item($r, 8).Value  = "Frutos de hueso (carozo)"
    $ws.Cells.Item($r, 9).Value  = 100103001
    $ws.Cells.Item($r, 10).Value = "Cereza"
    $ws.Cells.Item($r, 11).Value = "Rainier"
    $ws.Cells.Item($r, 12).Value = $quality[$i]
    $ws.Cells.Item($r, 13).Value = $volume[$i]
    $ws.Cells.Item($r, 14).Value = $minPrice[$i]
    $ws.Cells.Item($r, 15).Value = $maxPrice[$i]
    $ws.Cells.Item($r, 16).Value = $avgPrice[$i]
    $ws.Cells.Item($r, 17).Value = "`$/bandeja 10 kilos"
    $ws.Cells.Item($r, 18).Value = "Provincia de Curicó"
    $ws.Cells.Item($r, 19).Value = $perKg[$i]
    $ws.Cells.Item($r, 20).Value = 10
}
